# Update Asura_Profits market-data snapshot values (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 103.41304
$ws.Range("I33").Value = 96.19512
$ws.Range("J33").Value = 162.6
$ws.Range("K33").Value = 96.19512
$ws.Range("L33").Value = 162.6
$ws.Range("M33").Value = 132.80488
$ws.Range("N33").Value = -620.6
$ws.Range("H86").Value = 1968.5
$ws.Range("I86").Value = 2006.8572
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 2006.8572
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -883.8571999999999
$ws.Range("N86").Value = -3946
$ws.Range("H89").Value = 1968.5
$ws.Range("I89").Value = 2006.8572
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 10034.286
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -4418.286
$ws.Range("N89").Value = -19732
$ws.Range("H107").Value = 48126.906
$ws.Range("I107").Value = 62819.625
$ws.Range("J107").Value = 1110.2
$ws.Range("K107").Value = 62819.625
$ws.Range("L107").Value = 1110.2
$ws.Range("M107").Value = -60899.625
$ws.Range("N107").Value = -4950.2
$ws.Range("H138").Value = 3041.6296
$ws.Range("I138").Value = 1399.3334
$ws.Range("J138").Value = 5833.533
$ws.Range("K138").Value = 4198.0002
$ws.Range("L138").Value = 17500.599
$ws.Range("M138").Value = 941.9997999999996
$ws.Range("N138").Value = -27780.599
$ws.Range("H141").Value = 3907.6863
$ws.Range("I141").Value = 1757.4
$ws.Range("J141").Value = 20034.834
$ws.Range("K141").Value = 5272.200000000001
$ws.Range("L141").Value = 60104.50199999999
$ws.Range("M141").Value = -92.20000000000073
$ws.Range("N141").Value = -70464.50199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30345.902
$ws.Range("I32").Value = 9497.968999999999
$ws.Range("J32").Value = 197129.38
$ws.Range("K32").Value = 9497.968999999999
$ws.Range("L32").Value = 197129.38
$ws.Range("M32").Value = -9210.968999999999
$ws.Range("N32").Value = -197703.38
$ws.Range("H96").Value = 180114.67
$ws.Range("J96").Value = 180114.67
$ws.Range("L96").Value = 180114.67
$ws.Range("N96").Value = -185606.67
$ws.Range("H125").Value = 70714.86
$ws.Range("J125").Value = 70714.86
$ws.Range("L125").Value = 70714.86
$ws.Range("N125").Value = -80554.86
$ws.Range("H132").Value = 1763.0139
$ws.Range("I132").Value = 1336.3962
$ws.Range("J132").Value = 2953.0527
$ws.Range("K132").Value = 4009.188599999999
$ws.Range("L132").Value = 8859.158100000001
$ws.Range("M132").Value = -1479.188599999999
$ws.Range("N132").Value = -13919.1581

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 34634.676
$ws.Range("I20").Value = 48985.953
$ws.Range("J20").Value = 4497
$ws.Range("K20").Value = 48985.953
$ws.Range("L20").Value = 4497
$ws.Range("M20").Value = -48738.953
$ws.Range("N20").Value = -4991
$ws.Range("H86").Value = 205700.2
$ws.Range("I86").Value = 7333.3335
$ws.Range("J86").Value = 503250.5
$ws.Range("K86").Value = 7333.3335
$ws.Range("L86").Value = 503250.5
$ws.Range("M86").Value = -6210.3335
$ws.Range("N86").Value = -505496.5
$ws.Range("H89").Value = 205700.2
$ws.Range("I89").Value = 7333.3335
$ws.Range("J89").Value = 503250.5
$ws.Range("K89").Value = 36666.6675
$ws.Range("L89").Value = 2516252.5
$ws.Range("M89").Value = -31050.6675
$ws.Range("N89").Value = -2527484.5
$ws.Range("H100").Value = 30214
$ws.Range("J100").Value = 30214
$ws.Range("L100").Value = 30214
$ws.Range("N100").Value = -32378
$ws.Range("H109").Value = 21742.75
$ws.Range("J109").Value = 21742.75
$ws.Range("L109").Value = 21742.75
$ws.Range("N109").Value = -24516.75
$ws.Range("H124").Value = 24500
$ws.Range("J124").Value = 24500
$ws.Range("L124").Value = 24500
$ws.Range("N124").Value = -34320
$ws.Range("H134").Value = 2165.5095
$ws.Range("I134").Value = 1839.1052
$ws.Range("J134").Value = 2992.4
$ws.Range("K134").Value = 5517.3156
$ws.Range("L134").Value = 8977.200000000001
$ws.Range("M134").Value = -2982.3156
$ws.Range("N134").Value = -14047.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9000
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 50000
$ws.Range("K23").Value = 50000
$ws.Range("M23").Value = -49760
$ws.Range("H27").Value = 50000
$ws.Range("I27").Value = 50000
$ws.Range("K27").Value = 50000
$ws.Range("M27").Value = -49808
$ws.Range("H31").Value = 1700.5094
$ws.Range("I31").Value = 1342.9762
$ws.Range("J31").Value = 3065.6365
$ws.Range("K31").Value = 1342.9762
$ws.Range("L31").Value = 3065.6365
$ws.Range("M31").Value = -1047.9762
$ws.Range("N31").Value = -3655.6365
$ws.Range("H34").Value = 1700.5094
$ws.Range("I34").Value = 1342.9762
$ws.Range("J34").Value = 3065.6365
$ws.Range("K34").Value = 1342.9762
$ws.Range("L34").Value = 3065.6365
$ws.Range("M34").Value = -1140.9762
$ws.Range("N34").Value = -3469.6365
$ws.Range("H132").Value = 288672.7
$ws.Range("I132").Value = 330526.9
$ws.Range("J132").Value = 2668.8333
$ws.Range("K132").Value = 991580.7000000001
$ws.Range("L132").Value = 8006.499899999999
$ws.Range("M132").Value = -989050.7000000001
$ws.Range("N132").Value = -13066.4999
$ws.Range("H134").Value = 1080.0526
$ws.Range("I134").Value = 852.8125
$ws.Range("J134").Value = 2292
$ws.Range("K134").Value = 2558.4375
$ws.Range("L134").Value = 6876
$ws.Range("M134").Value = -23.4375
$ws.Range("N134").Value = -11946

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 523.4545000000001
$ws.Range("I113").Value = 482.52942
$ws.Range("J113").Value = 566.9375
$ws.Range("K113").Value = 1447.58826
$ws.Range("L113").Value = 1700.8125
$ws.Range("M113").Value = 722.41174
$ws.Range("N113").Value = -6040.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6468.6787
$ws.Range("I70").Value = 6107.1875
$ws.Range("J70").Value = 6950.6665
$ws.Range("K70").Value = 6107.1875
$ws.Range("L70").Value = 6950.6665
$ws.Range("M70").Value = -5837.1875
$ws.Range("N70").Value = -7490.6665
$ws.Range("H73").Value = 6468.6787
$ws.Range("I73").Value = 6107.1875
$ws.Range("J73").Value = 6950.6665
$ws.Range("K73").Value = 6107.1875
$ws.Range("L73").Value = 6950.6665
$ws.Range("M73").Value = -5171.1875
$ws.Range("N73").Value = -8822.666499999999
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("H102").Value = 3253.3914
$ws.Range("I102").Value = 2920
$ws.Range("J102").Value = 3772
$ws.Range("K102").Value = 2920
$ws.Range("L102").Value = 3772
$ws.Range("M102").Value = -1298
$ws.Range("N102").Value = -7016
$ws.Range("H123").Value = 8625
$ws.Range("J123").Value = 8625
$ws.Range("L123").Value = 8625
$ws.Range("N123").Value = -13525
$ws.Range("H132").Value = 1722.6522
$ws.Range("I132").Value = 1172.8064
$ws.Range("J132").Value = 2859
$ws.Range("K132").Value = 3518.4192
$ws.Range("L132").Value = 8577
$ws.Range("M132").Value = -988.4191999999998
$ws.Range("N132").Value = -13637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6000
$ws.Range("I2").Value = 6000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -5888
$ws.Range("H31").Value = 4114.4443
$ws.Range("I31").Value = 2028.5
$ws.Range("J31").Value = 5783.2
$ws.Range("K31").Value = 2028.5
$ws.Range("L31").Value = 5783.2
$ws.Range("M31").Value = -1780.5
$ws.Range("N31").Value = -6279.2
$ws.Range("H136").Value = 1435.211
$ws.Range("I136").Value = 1167.8987
$ws.Range("J136").Value = 3355
$ws.Range("K136").Value = 3503.6961
$ws.Range("L136").Value = 10065
$ws.Range("M136").Value = -953.6961000000001
$ws.Range("N136").Value = -15165

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 50000
$ws.Range("I11").Value = 50000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 50000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -49858
$ws.Range("H16").Value = 36959
$ws.Range("J16").Value = 36959
$ws.Range("L16").Value = 36959
$ws.Range("N16").Value = -37543
$ws.Range("H123").Value = 22227
$ws.Range("J123").Value = 22227
$ws.Range("L123").Value = 22227
$ws.Range("N123").Value = -32027
$ws.Range("H136").Value = 118755.68
$ws.Range("I136").Value = 1170.5807
$ws.Range("J136").Value = 1333801.6
$ws.Range("K136").Value = 3511.7421
$ws.Range("L136").Value = 4001404.8
$ws.Range("M136").Value = -961.7420999999999
$ws.Range("N136").Value = -4006504.8

# Columns that are fully cleared (no longer applicable) for a few rows
$wb.Worksheets.Item("GSM").Range("N98").ClearContents()
$wb.Worksheets.Item("LTW").Range("N2").ClearContents()
$wb.Worksheets.Item("WVR").Range("N11").ClearContents()
